$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Journal URL:" paragraph - replace the placeholder text with the actual
#    GitHub repository URL, keeping "Journal URL:" bold and the new text in
#    its own (non-bold) runs.
# ---------------------------------------------------------------------------
$journalPos = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Journal URL:*") {
        $journalPos = $i
        break
    }
}

if ($journalPos -gt 0) {
    $journalPara = $d.Paragraphs.Item($journalPos)
    $paraRange = $journalPara.Range

    # Remove everything after "Journal URL:" (the old placeholder text,
    # including the "[Insert Publicly-accessible Cloud Service URL]" part).
    $findRange = $d.Range($paraRange.Start, $paraRange.End)
    $findRange.Find.Execute(" [Insert Publicly-accessible Cloud Service URL]", `
        $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

    # Re-fetch the (now shorter) paragraph and append a plain space run.
    $journalPara = $d.Paragraphs.Item($journalPos)
    $journalPara.Range.InsertAfter(" ")

    # Make sure the space we just inserted is explicitly not bold (forces an
    # explicit - if empty - run-properties element, splitting it from the
    # bold "Journal URL:" run).
    $journalPara = $d.Paragraphs.Item($journalPos)
    $spaceStart = $journalPara.Range.End - 2
    $spaceRange = $d.Range($spaceStart, $spaceStart + 1)
    $spaceRange.Bold = 0

    # Append the URL text as its own run.
    $urlText = "https://github.com/prachijpatel/Learning_Journal"
    $journalPara = $d.Paragraphs.Item($journalPos)
    $urlInsertStart = $journalPara.Range.End - 1
    $journalPara.Range.InsertAfter($urlText)

    # Explicitly force not-bold on the new URL run too, so it gets its own
    # (empty) run-properties element, matching the target structure.
    $urlRange = $d.Range($urlInsertStart, $urlInsertStart + $urlText.Length)
    $urlRange.Bold = 1
    $urlRange.Bold = 0
}

# ---------------------------------------------------------------------------
# 2. "Personal Growth" reflection paragraph - switch it to single line
#    spacing (adds w:line="240" w:lineRule="auto").
# ---------------------------------------------------------------------------
$growthPos = 0
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Through this course,*solid foundation for future courses and projects.*") {
        $growthPos = $i
        break
    }
}

if ($growthPos -gt 0) {
    $growthPara = $d.Paragraphs.Item($growthPos)
    $growthPara.Format.LineSpacingRule = 0
}

# ---------------------------------------------------------------------------
# 3. Remove the trailing empty paragraphs at the very end of the document
#    (right before the sectPr), added after the "Personal Growth" reflection.
# ---------------------------------------------------------------------------
if ($growthPos -gt 0) {
    $count = $d.Paragraphs.Count
    if ($count -gt $growthPos) {
        $firstTrailing = $d.Paragraphs.Item($growthPos + 1)
        $lastPara = $d.Paragraphs.Item($count)
        $trailingRange = $d.Range($firstTrailing.Range.Start, $lastPara.Range.End)
        $trailingRange.Delete()
    }
}
